# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" worksheet (inserted between "总计" and "2022-Q2"),
# fills it with the Q3 fund-holdings data, and updates the "总计" (totals)
# sheet so that it lists the new Q3 entry on top of the existing Q2 / Q1 rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: push existing rows down one slot and
#    insert the new 2022-Q3 summary row at the top of the data (row 2).
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Duplicate formatting downward so the new / shifted rows keep the same look
# (bold, centered, bordered index column) as the existing data rows.
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A4:D4").PasteSpecial(-4122)

$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4122)

# Old row 3 (2022-Q1) moves to row 4
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(4,2).Value = "2022-Q1"
$wsTotal.Cells.Item(4,3).Value = 1
$wsTotal.Cells.Item(4,4).Value = 0.78

# Old row 2 (2022-Q2) moves to row 3
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2022-Q2"
$wsTotal.Cells.Item(3,3).Value = 16
$wsTotal.Cells.Item(3,4).Value = 3.34

# New row 2 (2022-Q3)
$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 13
$wsTotal.Cells.Item(2,4).Value = 2.11

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet and place it right before "2022-Q2".
# ---------------------------------------------------------------------------
$wsQ2  = $wb.Worksheets.Item(2)
$wsQ3  = $wb.Worksheets.Add()
$wsQ3.Name = "2022-Q3"

# Copy the header styling (bold / border / centered) from the Q2 sheet.
$wsQ2.Range("B1:H1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

# Copy the index-column styling (A2:A14) from the Q2 sheet's A column.
$wsQ2.Range("A2:A14").Copy()
$wsQ3.Range("A2:A14").PasteSpecial(-4122)

# Header row
$wsQ3.Cells.Item(1,2).Value = "基金代码"
$wsQ3.Cells.Item(1,3).Value = "基金名称"
$wsQ3.Cells.Item(1,4).Value = "基金规模"
$wsQ3.Cells.Item(1,5).Value = "股票总仓位"
$wsQ3.Cells.Item(1,6).Value = "仓位占比"
$wsQ3.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ3.Cells.Item(1,8).Value = "仓位排名"

# Columns B, D:G hold numeric-looking text values (fund codes / percentages)
# -- force text formatting first so things like "002345" or "0.1140" keep
# their leading / trailing zeros instead of being coerced into numbers.
# (Column C is always a non-numeric fund name, so it needs no special format.)
$wsQ3.Range("B2:B14").NumberFormat = "@"
$wsQ3.Range("D2:G14").NumberFormat = "@"

$q3Rows = @(
    @(0, "002345", "华夏高端制造灵活配置混合A",          "20.58", "90.97", "6.18", "1.2718", 3),
    @(1, "004640", "华夏节能环保股票A",                  "6.04",  "93.11", "4.74", "0.2863", 6),
    @(2, "012703", "华夏核心成长混合A",                  "4.94",  "93.34", "4.82", "0.2381", 7),
    @(3, "015229", "华夏低碳经济一年持有混合A",          "2.49",  "91.97", "4.58", "0.1140", 7),
    @(4, "015058", "华夏高端制造灵活配置混合C",          "0.95",  "90.97", "6.18", "0.0587", 3),
    @(5, "015230", "华夏低碳经济一年持有混合C",          "1.02",  "91.97", "4.58", "0.0467", 7),
    @(6, "012710", "华夏核心成长混合C",                  "0.61",  "93.34", "4.82", "0.0294", 7),
    @(7, "002071", "长安产业精选灵活配置混合C",          "0.95",  "77.85", "3.01", "0.0286", 8),
    @(8, "000496", "长安产业精选灵活配置混合A",          "0.44",  "77.85", "3.01", "0.0132", 8),
    @(9, "015060", "华夏节能环保股票C",                  "0.26",  "93.11", "4.74", "0.0123", 6),
    @(10,"004536", "嘉实中小企业量化活力灵活配置混合",    "0.21",  "94.75", "4.03", "0.0085", 7),
    @(11,"001281", "长安鑫利优选灵活配置混合A",          "0.13",  "80.57", "3.08", "0.0040", 9),
    @(12,"002072", "长安鑫利优选灵活配置混合C",          "0.10",  "80.57", "3.08", "0.0031", 9)
)

$r = 2
foreach ($row in $q3Rows) {
    $wsQ3.Cells.Item($r,1).Value = $row[0]
    $wsQ3.Cells.Item($r,2).Value = $row[1]
    $wsQ3.Cells.Item($r,3).Value = $row[2]
    $wsQ3.Cells.Item($r,4).Value = $row[3]
    $wsQ3.Cells.Item($r,5).Value = $row[4]
    $wsQ3.Cells.Item($r,6).Value = $row[5]
    $wsQ3.Cells.Item($r,7).Value = $row[6]
    $wsQ3.Cells.Item($r,8).Value = $row[7]
    $r++
}

# Move the freshly built sheet so the final order is:
# 总计, 2022-Q3, 2022-Q2, 2022-Q1
$wsQ3.Move($wsQ2)
